# Apply "Generate Report for Handoff" update:
#  - file b1f45bce-ead6-4194-8875-91118ab5a876.md is renamed/replaced by
#    3f31608e-0d5b-405d-b49e-acdd59fef6b8.md
#  - file b2fd5750-ae99-46e5-8811-4c47ddc44e2f.md is renamed/replaced by
#    ffffa091a720-58b1-49a1-a246-fdf430061668.md
#  - status changes from "Handed back: in sync with en-US" to "Ready for handoff"
#  - handoff xliff + datetimes refreshed, both language rows now point at the
#    SAME freshly generated handoff xliff (content-duplicate), handback info
#    cleared out (no handback has happened yet for this new handoff round).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "3f31608e-0d5b-405d-b49e-acdd59fef6b8.md"
$ov.Range("C2").Value = ".md"
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-26 07:01:21"

$ov.Range("A3").Value = "ffffa091a720-58b1-49a1-a246-fdf430061668.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-26 07:01:21"

# Recreate the two hyperlinks (B2/B3) with the new display text, pointing at
# the same targets the old hyperlinks used.
$ovUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/b1f45bce-ead6-4194-8875-91118ab5a876.md"
$ovUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/b2fd5750-ae99-46e5-8811-4c47ddc44e2f.md"
$ov.Hyperlinks.Delete()
$ov.Range("B2").Value = "e2e\3f31608e-0d5b-405d-b49e-acdd59fef6b8.md"
$ov.Range("B3").Value = "e2e\ffffa091a720-58b1-49a1-a246-fdf430061668.md"
[void]$ov.Hyperlinks.Add($ov.Range("B2"), $ovUrl1, "", "", "e2e\3f31608e-0d5b-405d-b49e-acdd59fef6b8.md")
[void]$ov.Hyperlinks.Add($ov.Range("B3"), $ovUrl2, "", "", "e2e\ffffa091a720-58b1-49a1-a246-fdf430061668.md")

$ov.Columns.Item(5).AutoFit()
$ov.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "3f31608e-0d5b-405d-b49e-acdd59fef6b8.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("F2").Value = "False"
$zh.Range("G2").Value = "3f31608e-0d5b-405d-b49e-acdd59fef6b8.ae2eb6082684251e7dc10f7ff4183c6dc6d4e931.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-26 07:01:16"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = "ffffa091a720-58b1-49a1-a246-fdf430061668.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = "3f31608e-0d5b-405d-b49e-acdd59fef6b8.ae2eb6082684251e7dc10f7ff4183c6dc6d4e931.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-26 07:01:16"
$zh.Range("I3").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"

$zhUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/b1f45bce-ead6-4194-8875-91118ab5a876.md"
$zhUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/b2fd5750-ae99-46e5-8811-4c47ddc44e2f.md"
$zh.Hyperlinks.Delete()
[void]$zh.Hyperlinks.Add($zh.Range("A2"), $zhUrl1, "", "", "3f31608e-0d5b-405d-b49e-acdd59fef6b8.md")
[void]$zh.Hyperlinks.Add($zh.Range("A3"), $zhUrl2, "", "", "ffffa091a720-58b1-49a1-a246-fdf430061668.md")
# I2/I3 lost their hyperlink style (normal body text now, cell is blank).
$zh.Range("I2").Style = "Normal"
$zh.Range("I3").Style = "Normal"

$zh.Columns.Item(3).AutoFit()
$zh.Columns.Item(9).AutoFit()
$zh.Columns.Item(10).AutoFit()

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "3f31608e-0d5b-405d-b49e-acdd59fef6b8.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("F2").Value = "False"
$de.Range("G2").Value = "3f31608e-0d5b-405d-b49e-acdd59fef6b8.ae2eb6082684251e7dc10f7ff4183c6dc6d4e931.de-de.xlf"
$de.Range("H2").Value = "2016-08-26 07:01:21"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Range("A3").Value = "ffffa091a720-58b1-49a1-a246-fdf430061668.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "True"
$de.Range("G3").Value = "3f31608e-0d5b-405d-b49e-acdd59fef6b8.ae2eb6082684251e7dc10f7ff4183c6dc6d4e931.de-de.xlf"
$de.Range("H3").Value = "2016-08-26 07:01:21"
$de.Range("I3").Value = ""
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"

$deUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/b1f45bce-ead6-4194-8875-91118ab5a876.md"
$deUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/b2fd5750-ae99-46e5-8811-4c47ddc44e2f.md"
$de.Hyperlinks.Delete()
[void]$de.Hyperlinks.Add($de.Range("A2"), $deUrl1, "", "", "3f31608e-0d5b-405d-b49e-acdd59fef6b8.md")
[void]$de.Hyperlinks.Add($de.Range("A3"), $deUrl2, "", "", "ffffa091a720-58b1-49a1-a246-fdf430061668.md")
$de.Range("I2").Style = "Normal"
$de.Range("I3").Style = "Normal"

$de.Columns.Item(3).AutoFit()
$de.Columns.Item(9).AutoFit()
$de.Columns.Item(10).AutoFit()
